$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7512820512820513
$ws.Range("D2").Value = 0.6563714620297117
$ws.Range("E2").Value = 0.5559950813075812
$ws.Range("F2").Value = 0.7096774193548386
$ws.Range("G2").Value = 0.6441498937684872
$ws.Range("H2").Value = 0.5479141865079364
$ws.Range("I2").Value = 0.7333333333333333
$ws.Range("J2").Value = 0.6924867021276595
$ws.Range("K2").Value = 0.5986979166666666
$ws.Range("L2").Value = 0.6875
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 0 0 1 1 1 0 0 1 0 0 1 1 1 0 0 1 1 1 1 1 1 1 0]'
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()), (''selector'', None),
                (''model'',
                 LogisticRegression(C=0.001, max_iter=1000, penalty=''l1'',
                                    random_state=42, solver=''saga''))])'
$ws.Range("B3").Value = 0.7499999999999999
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': RobustScaler(), ''model__solver'': ''saga'', ''model__penalty'': ''l1'', ''model__class_weight'': None, ''model__C'': 0.001}'
$ws.Range("D3").Value = 0.6650555282310263
$ws.Range("E3").Value = 0.5591489801048625
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 0.6368700552463364
$ws.Range("H3").Value = 0.5822584033613446
$ws.Range("I3").Value = 0.6666666666666666
$ws.Range("J3").Value = 0.7230600750938673
$ws.Range("K3").Value = 0.5838235294117646
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa0013b550>),
                (''model'',
                 LogisticRegression(C=5, max_iter=1000, random_state=42,
                                    solver=''liblinear''))])'
$ws.Range("B4").Value = 0.6941025641025641
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa00129a90>, ''scaler'': MinMaxScaler(), ''model__solver'': ''liblinear'', ''model__penalty'': ''l2'', ''model__class_weight'': None, ''model__C'': 5}'
$ws.Range("D4").Value = 0.6351495788518212
$ws.Range("E4").Value = 0.539090197996448
$ws.Range("F4").Value = 0.7058823529411765
$ws.Range("G4").Value = 0.6217265611945623
$ws.Range("H4").Value = 0.5364186507936508
$ws.Range("I4").Value = 0.8
$ws.Range("J4").Value = 0.6663194444444444
$ws.Range("K4").Value = 0.5840625
$ws.Range("L4").Value = 0.631578947368421
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[1 1 1 0 1 1 0 1 1 0 1 0 0 1 1 1 0 1 0 0 1 0 1 1]'
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa00129640>),
                (''model'',
                 LogisticRegression(C=0.0001, max_iter=1000, random_state=42,
                                    solver=''saga''))])'
$ws.Range("B5").Value = 0.7664285714285713
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa00039f10>, ''scaler'': None, ''model__solver'': ''saga'', ''model__penalty'': ''l2'', ''model__class_weight'': None, ''model__C'': 0.0001}'
$ws.Range("D5").Value = 0.6633817260601386
$ws.Range("E5").Value = 0.5895066218503718
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.635783598306968
$ws.Range("H5").Value = 0.5703980654761904
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.7242346938775509
$ws.Range("K5").Value = 0.6562500000000001
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', None),
                (''model'',
                 LogisticRegression(C=0.0001, max_iter=1000, random_state=42,
                                    solver=''liblinear''))])'
$ws.Range("B6").Value = 0.7499999999999999
$ws.Range("C6").Value = '{''selector'': None, ''scaler'': MinMaxScaler(), ''model__solver'': ''liblinear'', ''model__penalty'': ''l2'', ''model__class_weight'': None, ''model__C'': 0.0001}'
$ws.Range("D6").Value = 0.6306959182515954
$ws.Range("E6").Value = 0.5261904034035181
$ws.Range("F6").Value = 0.6285714285714286
$ws.Range("G6").Value = 0.6044979785268462
$ws.Range("H6").Value = 0.5028011969815248
$ws.Range("I6").Value = 0.4583333333333333
$ws.Range("J6").Value = 0.6764501891551072
$ws.Range("K6").Value = 0.5816939890710383
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
